$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A21").Value = "Justine Calma"
$ws.Range("B21").Value = "Tesla, led by Elon Musk, confirmed that it purchased about `$ 1.5 billion in bitcoin in January and expects to start accepting it as a payment in the future.  | Photo by Artur Widak/NurPhoto via Getty… [+6118 chars]"
$ws.Range("C21").Value = "Tesla announced this week that it purchased `$1.5 billion in bitcoin and plans to accept bitcoin as payment in the future. CEO Elon Musk might have boosted bitcoin’s credibility. But greenhouse gas emissions from bitcoin are on the rise as the price of the cry…"
$ws.Range("D21").Value = "2021-02-09T23:10:46Z"
$ws.Range("E21").Value = "{'id': 'the-verge', 'name': 'The Verge'}"
$ws.Range("F21").Value = "Tesla’s `$1.5 billion bitcoin purchase clashes with its environmental aspirations"
$ws.Range("G21").Value = "https://www.theverge.com/2021/2/9/22275243/teslas-bitcoin-purchase-clashes-climate-change-mission"
$ws.Range("H21").Value = "https://cdn.vox-cdn.com/thumbor/eQST6cDSwG_MASzsxTUv3igFdxU=/0x430:6000x3571/fit-in/1200x630/cdn.vox-cdn.com/uploads/chorus_asset/file/22293241/1231048263.jpg"
$ws.Range("I21").Value = "2021.02.09"
$ws.Range("J21").Value = 0

$ws.Range("A22").Value = "Gregory Barber"
$ws.Range("B22").Value = "For a brief moment on Sunday, before Tesla said it had invested `$1.5 billion in bitcoin and planned to let people use the cryptocurrency to pay for its cars, bitcoins price could be expressed with a … [+3302 chars]"
$ws.Range("C22").Value = "The maker of electric vehicles said it had invested `$1.5 billion in bitcoin and plans to accept the cryptocurrency as payment for its cars."
$ws.Range("D22").Value = "2021-02-09T12:00:00Z"
$ws.Range("E22").Value = "{'id': 'wired', 'name': 'Wired'}"
$ws.Range("F22").Value = "Would You Trade a Bitcoin for a Tesla?"
$ws.Range("G22").Value = "https://www.wired.com/story/would-you-trade-bitcoin-tesla/"
$ws.Range("H22").Value = "https://media.wired.com/photos/6021da4ff1bf194f33695dc3/191:100/w_1280,c_limit/business_tesla-bitcoin_1229893729.jpg"
$ws.Range("I22").Value = "2021.02.09"
$ws.Range("J22").Value = 0

$ws.Range("A23").Value = "Manish Singh"
$ws.Range("B23").Value = "Twitter and Square CEO Jack Dorsey and rapper Jay Z have created an endowment to fund bitcoin development initially in Africa and India, Dorsey said Friday.`nThe duo is putting 500 bitcoin, which is … [+3984 chars]"
$ws.Range("C23").Value = "Twitter and Square CEO Jack Dorsey and rapper Jay Z have created an endowment to fund bitcoin development initially in Africa and India, Dorsey said Friday. The duo is putting 500 bitcoin, which is currently worth `$23.6 million, in the endowment called ₿trust…"
$ws.Range("D23").Value = "2021-02-12T10:51:25Z"
$ws.Range("E23").Value = "{'id': 'techcrunch', 'name': 'TechCrunch'}"
$ws.Range("F23").Value = "Jack Dorsey and Jay Z invest `$23.6 million to fund Bitcoin development"
$ws.Range("G23").Value = "http://techcrunch.com/2021/02/12/jack-dorsey-and-jay-z-invest-23-6-million-to-fund-bitcoin-development/"
$ws.Range("H23").Value = "https://techcrunch.com/wp-content/uploads/2020/11/GettyImages-887657568.jpg?w=600"
$ws.Range("I23").Value = "2021.02.12"
$ws.Range("J23").Value = 0.15

$ws.Range("A24").Value = "Alex Wilhelm"
$ws.Range("B24").Value = "Today in an SEC filing, Tesla disclosed that it has acquired `$1.5 billion in bitcoin, the popular cryptocurrency. Moreover, the company noted that it may also accept bitcoin in the future as a form o… [+2230 chars]"
$ws.Range("C24").Value = "Today in an SEC filing, Tesla disclosed that it has acquired `$1.5 billion in bitcoin, the popular cryptocurrency. Moreover, the company noted that it may also accept bitcoin in the future as a form of payment for its cars, though it did allow that there is so…"
$ws.Range("D24").Value = "2021-02-08T13:02:03Z"
$ws.Range("E24").Value = "{'id': 'techcrunch', 'name': 'TechCrunch'}"
$ws.Range("F24").Value = "Tesla buys `$1.5B in bitcoin, may accept the cryptocurrency as payment in the future"
$ws.Range("G24").Value = "http://techcrunch.com/2021/02/08/tesla-buys-1-5b-in-bitcoin-may-accept-the-cryptocurrency-as-payment-in-the-future/"
$ws.Range("H24").Value = "https://techcrunch.com/wp-content/uploads/2020/09/tesla-glitch1.jpg?w=712"
$ws.Range("I24").Value = "2021.02.08"
$ws.Range("J24").Value = 0.1666666666666667

$ws.Range("A25").Value = "WIRED Staff"
$ws.Range("B25").Value = "When Bitcoin first appeared out of digital thin air, it was hailed as having the potential to upend the way people spent money. But more than a decade later, cryptocurrency is still only trickling in… [+2039 chars]"
$ws.Range("C25").Value = "This week, we discuss crypto’s role in the future of shopping. When can we to use it to buy everything from Nikes to Teslas?"
$ws.Range("D25").Value = "2021-02-12T13:00:00Z"
$ws.Range("E25").Value = "{'id': 'wired', 'name': 'Wired'}"
$ws.Range("F25").Value = "This Cryptocurrency Is Really Burning a Hole in My Pocket"
$ws.Range("G25").Value = "https://www.wired.com/story/gadget-lab-podcast-491/"
$ws.Range("H25").Value = "https://media.wired.com/photos/602589a9e27a393fd0c185af/191:100/w_1280,c_limit/Gear-Tesla-GL-Podcast-1227575736.jpg"
$ws.Range("I25").Value = "2021.02.12"
$ws.Range("J25").Value = 0

$ws.Range("A26").Value = "Anthony Ha"
$ws.Range("B26").Value = "Jack Dorsey and Jay Z create a bitcoin endowment, Datadog acquires a Startup Battlefield company and BuzzFeed experiments with AI-generated quizzes. This is your Daily Crunch for February 12, 2021.`n… [+2805 chars]"
$ws.Range("C26").Value = "Jack Dorsey and Jay Z create a bitcoin endowment, Datadog acquires a Startup Battlefield company and BuzzFeed experiments with AI-generated quizzes. This is your Daily Crunch for February 12, 2021. Oh, and before we get started: Consider applying to the Early…"
$ws.Range("D26").Value = "2021-02-12T23:10:51Z"
$ws.Range("E26").Value = "{'id': 'techcrunch', 'name': 'TechCrunch'}"
$ws.Range("F26").Value = "Daily Crunch: Jack Dorsey and Jay Z invest in bitcoin development"
$ws.Range("G26").Value = "http://techcrunch.com/2021/02/12/daily-crunch-jack-dorsey-and-jay-z-invest-in-bitcoin-development/"
$ws.Range("H26").Value = "https://techcrunch.com/wp-content/uploads/2020/07/GettyImages-887657568.jpg?w=600"
$ws.Range("I26").Value = "2021.02.12"
$ws.Range("J26").Value = 0

$ws.Range("A27").Value = "Lucas Matney"
$ws.Range("B27").Value = "Cryptocurrencies, more so than most other things, are only valuable because of a shared agreement that they are valuable. Their value is a product of digital handshakes over millions of transactions … [+3710 chars]"
$ws.Range("C27").Value = "Cryptocurrencies, more so than most other things, are only valuable because of a shared agreement that they are valuable. Their value is a product of digital handshakes over millions of transactions firming up that consensus. For bitcoin, the trust that it ha…"
$ws.Range("D27").Value = "2021-01-27T19:24:32Z"
$ws.Range("E27").Value = "{'id': 'techcrunch', 'name': 'TechCrunch'}"
$ws.Range("F27").Value = "Could meme stocks like GameStop kill bitcoin’s rise?"
$ws.Range("G27").Value = "http://techcrunch.com/2021/01/27/could-meme-stocks-like-gamestop-kill-bitcoins-rise/"
$ws.Range("H27").Value = "https://techcrunch.com/wp-content/uploads/2017/08/bitcoin-split-2017a.jpg?w=711"
$ws.Range("I27").Value = "2021.01.27"
$ws.Range("J27").Value = 0.175

$ws.Range("A28").Value = "Lucas Matney"
$ws.Range("B28").Value = "Coinbase plans to go public by way of a direct listing, the company announced in a blog post today.`nThe cryptocurrency exchange was founded in 2012 and allows users to buy and trade decentralized to… [+1258 chars]"
$ws.Range("C28").Value = "Coinbase plans to go public by way of a direct listing, the company announced in a blog post today. The cryptocurrency exchange was founded in 2012 and allows users to buy and trade decentralized tokens like bitcoin and ethereum. The company has raised over `$…"
$ws.Range("D28").Value = "2021-01-28T19:36:35Z"
$ws.Range("E28").Value = "{'id': 'techcrunch', 'name': 'TechCrunch'}"
$ws.Range("F28").Value = "Coinbase is going public via direct listing"
$ws.Range("G28").Value = "http://techcrunch.com/2021/01/28/coinbase-is-going-public-via-direct-listing/"
$ws.Range("H28").Value = "https://techcrunch.com/wp-content/uploads/2019/01/coinbase.jpg?w=750"
$ws.Range("I28").Value = "2021.01.28"
$ws.Range("J28").Value = 0.05

$ws.Range("A29").Value = "Romain Dillet"
$ws.Range("B29").Value = "Blockchain.com has announced that it has raised a `$120 million funding round. The company develops a popular cryptocurrency wallet as well as an exchange, an explorer and more.`nMoore Strategic Ventu… [+2165 chars]"
$ws.Range("C29").Value = "Blockchain.com has announced that it has raised a `$120 million funding round. The company develops a popular cryptocurrency wallet as well as an exchange, an explorer and more. Moore Strategic Ventures, Kyle Bass, Access Industries, Rovida Advisors, Lightspee…"
$ws.Range("D29").Value = "2021-02-17T18:24:31Z"
$ws.Range("E29").Value = "{'id': 'techcrunch', 'name': 'TechCrunch'}"
$ws.Range("F29").Value = "Crypto wallet and exchange company Blockchain.com raises `$120 million"
$ws.Range("G29").Value = "http://techcrunch.com/2021/02/17/crypto-wallet-and-exchange-company-blockchain-com-raises-120-million/"
$ws.Range("H29").Value = "https://techcrunch.com/wp-content/uploads/2021/02/andre-francois-mckenzie-JrjhtBJ-pGU-unsplash.jpg?w=599"
$ws.Range("I29").Value = "2021.02.17"
$ws.Range("J29").Value = 0.1875

$ws.Range("A30").Value = "Manish Singh"
$ws.Range("B30").Value = "India plans to introduce a law to ban private cryptocurrencies such as bitcoin in the country and provide a framework for the creation of an official digital currency during the current budget sessio… [+1858 chars]"
$ws.Range("C30").Value = "India plans to introduce a law to ban private cryptocurrencies such as bitcoin in the country and provide a framework for the creation of an official digital currency during the current budget session of parliament. In the agenda (PDF) published on the lower …"
$ws.Range("D30").Value = "2021-01-30T11:31:02Z"
$ws.Range("E30").Value = "{'id': 'techcrunch', 'name': 'TechCrunch'}"
$ws.Range("F30").Value = "India plans to introduce law to ban Bitcoin, other private cryptocurrencies"
$ws.Range("G30").Value = "http://techcrunch.com/2021/01/30/india-plans-to-introduce-law-to-ban-bitcoin-other-private-cryptocurrencies/"
$ws.Range("H30").Value = "https://techcrunch.com/wp-content/uploads/2018/09/cryptocurrency.jpg?w=609"
$ws.Range("I30").Value = "2021.01.30"
$ws.Range("J30").Value = 0

$ws.Range("A31").Value = "Anthony Ha"
$ws.Range("B31").Value = "DoorDash acquires a salad-making robotics startup, Twitter confirms subscription plans and Tesla makes a big bet on bitcoin. This is your Daily Crunch for February 8, 2021.`nThe big story: DoorDash a… [+3188 chars]"
$ws.Range("C31").Value = "DoorDash acquires a salad-making robotics startup, Twitter confirms subscription plans and Tesla makes a big bet on bitcoin. This is your Daily Crunch for February 8, 2021. The big story: DoorDash acquires Chowbotics DoorDash has acquired the Bay Area startup…"
$ws.Range("D31").Value = "2021-02-08T23:27:07Z"
$ws.Range("E31").Value = "{'id': 'techcrunch', 'name': 'TechCrunch'}"
$ws.Range("F31").Value = "Daily Crunch: DoorDash acquires Chowbotics"
$ws.Range("G31").Value = "http://techcrunch.com/2021/02/08/daily-crunch-doordash-acquires-chowbotics/"
$ws.Range("H31").Value = "https://techcrunch.com/wp-content/uploads/2017/03/chowbotics_team_salads.jpg?w=600"
$ws.Range("I31").Value = "2021.02.08"
$ws.Range("J31").Value = 0

$ws.Range("A32").Value = "Reuters Staff"
$ws.Range("B32").Value = "By Reuters Staff`nSINGAPORE, Feb 9 (Reuters) - Cryptocurrencies extended gains in Asia on Tuesday, with bitcoin and ethereum reaching record highs, in the wake of a Tesla Inc investment in bitcoin.`n… [+389 chars]"
$ws.Range("C32").Value = "Cryptocurrencies extended gains in Asia on Tuesday, with bitcoin and ethereum reaching record highs, in the wake of a Tesla Inc investment in bitcoin."
$ws.Range("D32").Value = "2021-02-09T00:43:00Z"
$ws.Range("E32").Value = "{'id': 'reuters', 'name': 'Reuters'}"
$ws.Range("F32").Value = "Bitcoin extends gains above `$47,000 in Asia - Reuters"
$ws.Range("G32").Value = "https://www.reuters.com/article/crypto-currency-idUSL1N2KE33M"
$ws.Range("H32").Value = "https://s1.reutersmedia.net/resources_v2/images/rcom-default.png?w=800"
$ws.Range("I32").Value = "2021.02.09"
$ws.Range("J32").Value = 0

$ws.Range("A33").Value = "Reuters Staff"
$ws.Range("B33").Value = "By Reuters Staff`nFILE PHOTO: The logo of the Bitcoin digital currency is seen in a shop in Marseille, France, February 7, 2021. REUTERS/Eric Gaillard`nSINGAPORE (Reuters) - Cryptocurrencies extended… [+421 chars]"
$ws.Range("C33").Value = "Cryptocurrencies extended gains in Asia on Tuesday, with bitcoin and ethereum reaching record highs, in the wake of a Tesla Inc investment in bitcoin."
$ws.Range("D33").Value = "2021-02-09T00:48:00Z"
$ws.Range("E33").Value = "{'id': 'reuters', 'name': 'Reuters'}"
$ws.Range("F33").Value = "Bitcoin extends gains above `$47,000 in Asia - Reuters"
$ws.Range("G33").Value = "https://www.reuters.com/article/us-crypto-currency-idUSKBN2A902T"
$ws.Range("H33").Value = "https://static.reuters.com/resources/r/?m=02&d=20210209&t=2&i=1550815750&r=LYNXMPEH1801A&w=800"
$ws.Range("I33").Value = "2021.02.09"
$ws.Range("J33").Value = 0

$ws.Range("A34").Value = "Bhargav Acharya"
$ws.Range("B34").Value = "(Repeats for wider audience with no changes to text)`nFeb 9 (Reuters) - A Reddit user claiming to be a Tesla insider appeared to announce the carmakers purchase of bitcoin a month ago, according to a… [+1793 chars]"
$ws.Range("C34").Value = "A Reddit user claiming to be a Tesla insider appeared to announce the carmaker's purchase of bitcoin a month ago, according to a January post on the platform that said the electric carmaker had bought `$800 million worth of Bitcoin."
$ws.Range("D34").Value = "2021-02-09T17:34:00Z"
$ws.Range("E34").Value = "{'id': 'reuters', 'name': 'Reuters'}"
$ws.Range("F34").Value = "RPT-Reddit user claiming to be Tesla insider appeared to reveal bitcoin buy a month ago - Reuters"
$ws.Range("G34").Value = "https://www.reuters.com/article/crypto-currency-tesla-reddit-idUSL1N2KF29S"
$ws.Range("H34").Value = "https://s1.reutersmedia.net/resources_v2/images/rcom-default.png?w=800"
$ws.Range("I34").Value = "2021.02.09"
$ws.Range("J34").Value = 0.3

$ws.Range("A35").Value = "Bhargav Acharya"
$ws.Range("B35").Value = "FILE PHOTO: Representations of virtual currency Bitcoin are seen in front of Tesla logo in this illustration taken, February 9, 2021. REUTERS/Dado Ruvic/Illustration`n(Reuters) - A Reddit user claimi… [+1771 chars]"
$ws.Range("C35").Value = "A Reddit user claiming to be a Tesla insider appeared to announce the carmaker's purchase of bitcoin a month ago, according to a January post on the platform that said the electric carmaker had bought `$800 million worth of Bitcoin."
$ws.Range("D35").Value = "2021-02-09T15:16:00Z"
$ws.Range("E35").Value = "{'id': 'reuters', 'name': 'Reuters'}"
$ws.Range("F35").Value = "Reddit user claiming to be Tesla insider appeared to reveal bitcoin buy a month ago - Reuters UK"
$ws.Range("G35").Value = "https://www.reuters.com/article/crypto-currency-tesla-reddit-int-idUSKBN2A922X"
$ws.Range("H35").Value = "https://static.reuters.com/resources/r/?m=02&d=20210209&t=2&i=1550899933&r=LYNXMPEH181AF&w=800"
$ws.Range("I35").Value = "2021.02.09"
$ws.Range("J35").Value = 0.3

$ws.Range("A36").Value = "Gertrude Chavez-Dreyfuss"
$ws.Range("B36").Value = "* Dollar index little changed after Friday payrolls fall`n * Jobs data takes shine off dollar rebound`n * Ethereum gains on futures debut`n * Bitcoin hits record high after Tesla purchase`n * Graphic… [+4476 chars]"
$ws.Range("C36").Value = "* Dollar index little changed after Friday payrolls fall * Jobs data takes shine off dollar rebound * Ethereum gains on futures debut * Bitcoin hits record high after Tesla purchase * Graphic: World FX rates https://tmsnrt.rs/2RBWI5E (Adds details on Bitcoin,…"
$ws.Range("D36").Value = "2021-02-08T17:18:00Z"
$ws.Range("E36").Value = "{'id': 'reuters', 'name': 'Reuters'}"
$ws.Range("F36").Value = "FOREX-Dollar steadies after U.S. jobs-related losses - Reuters"
$ws.Range("G36").Value = "https://www.reuters.com/article/global-forex-idUSL1N2KE1TJ"
$ws.Range("H36").Value = "https://s1.reutersmedia.net/resources_v2/images/rcom-default.png?w=800"
$ws.Range("I36").Value = "2021.02.08"
$ws.Range("J36").Value = -0.009166666666666665

$ws.Range("A37").Value = "Reuters Staff"
$ws.Range("B37").Value = "By Reuters Staff`nFeb 21 (Reuters) - Bitcoin continued gaining on Sunday, rising to a fresh high and extending a two-month rally that took its market capitalization above `$1 trillion on Friday.`nThe … [+296 chars]"
$ws.Range("C37").Value = "Bitcoin continued gaining on Sunday, rising to a fresh high and extending a two-month rally that took its market capitalization above `$1 trillion on Friday."
$ws.Range("D37").Value = "2021-02-21T19:38:00Z"
$ws.Range("E37").Value = "{'id': 'reuters', 'name': 'Reuters'}"
$ws.Range("F37").Value = "Bitcoin hits fresh high - Reuters"
$ws.Range("G37").Value = "https://www.reuters.com/article/crypto-currency-bitcoin-idUSL1N2KR0FM"
$ws.Range("H37").Value = "https://s1.reutersmedia.net/resources_v2/images/rcom-default.png?w=800"
$ws.Range("I37").Value = "2021.02.21"
$ws.Range("J37").Value = 0.1533333333333333

$ws.Range("A38").Value = "Reuters Staff"
$ws.Range("B38").Value = "By Reuters Staff`nFILE PHOTO: Representations of virtual currency Bitcoin are seen in this picture illustration taken taken March 13, 2020. REUTERS/Dado Ruvic`nNEW YORK (Reuters) - Bitcoin fell sharp… [+254 chars]"
$ws.Range("C38").Value = "Bitcoin fell sharply on Wednesday, briefly slipping below the `$30,000 mark for the first time in five days, amid a broad risk-off move in global financial markets."
$ws.Range("D38").Value = "2021-01-27T15:55:00Z"
$ws.Range("E38").Value = "{'id': 'reuters', 'name': 'Reuters'}"
$ws.Range("F38").Value = "Bitcoin slips below `$30,000 - Reuters"
$ws.Range("G38").Value = "https://www.reuters.com/article/crypto-currencies-int-idUSKBN29W21S"
$ws.Range("H38").Value = "https://static.reuters.com/resources/r/?m=02&d=20210127&t=2&i=1549257422&r=LYNXMPEH0Q16C&w=800"
$ws.Range("I38").Value = "2021.01.27"
$ws.Range("J38").Value = 0.01458333333333333

$ws.Range("A39").Value = "Natasha Mascarenhas"
$ws.Range("B39").Value = "Hello and welcome back to Equity, TechCrunchs venture capital-focused podcast, where we unpack the numbers behind the headlines.`nNatasha and Danny and Alex and Grace were all here to chat through th… [+2406 chars]"
$ws.Range("C39").Value = "Hello and welcome back to Equity, TechCrunch’s venture capital-focused podcast, where we unpack the numbers behind the headlines. Natasha and Danny and Alex and Grace were all here to chat through the week’s biggest tech happenings. In very good Show News™, C…"
$ws.Range("D39").Value = "2021-02-18T22:37:32Z"
$ws.Range("E39").Value = "{'id': 'techcrunch', 'name': 'TechCrunch'}"
$ws.Range("F39").Value = "A16z doesn’t invest, it manifests"
$ws.Range("G39").Value = "http://techcrunch.com/2021/02/18/a16z-doesnt-invest-it-manifests/"
$ws.Range("H39").Value = "https://techcrunch.com/wp-content/uploads/2019/11/equity-podcast-2019-phone-1-3.jpg?w=430"
$ws.Range("I39").Value = "2021.02.18"
$ws.Range("J39").Value = 0.3275

$ws.Range("A40").Value = "Noel Randewich"
$ws.Range("B40").Value = "Feb 8 (Reuters) - Teslas `$1.5 billion Bitcoin investment unveiled on Monday may not surprise CEO Elon Musks followers on Twitter, where he has shown himself as a major proponent of the soaring crypto… [+1614 chars]"
$ws.Range("C40").Value = "Tesla's `$1.5 billion Bitcoin investment unveiled on Monday may not surprise CEO Elon Musk's followers on Twitter, where he has shown himself as a major proponent of the soaring cryptocurrency."
$ws.Range("D40").Value = "2021-02-08T20:58:00Z"
$ws.Range("E40").Value = "{'id': 'reuters', 'name': 'Reuters'}"
$ws.Range("F40").Value = "Musk's Bitcoin investment follows months of Twitter talk - Reuters"
$ws.Range("G40").Value = "https://www.reuters.com/article/crypto-currency-tesla-tweets-idUSL1N2KE263"
$ws.Range("H40").Value = "https://s1.reutersmedia.net/resources_v2/images/rcom-default.png?w=800"
$ws.Range("I40").Value = "2021.02.08"
$ws.Range("J40").Value = 0.0625

